$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table")

# Helper: assign a value as genuine Text (avoids Excel's automatic
# "looks like a number" coercion for strings such as "3", "10", "11"),
# while leaving the cell on the default/unstyled format afterwards.
function Set-TextValue {
    param($rng, $val)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Insert a new column in the "Table" sheet right before column E ---
# (Conf_index=B, 1=C, 2=D stay put; everything from the old "3" column
# onward shifts one slot to the right, and a brand-new "3" column lands
# in the freed-up E slot.)
$ws.Range("E1").EntireColumn.Insert()

# Grow the Tableau1 list-object so it covers the freshly inserted column
# plus the new trailing column (now B1:M7 instead of B1:L7).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B1:M7"))

# Row 2 (HW_version): rename the "system_1" configuration and add the
# new "system_1 (-O0)" configuration right next to it. Done before the
# header re-typing below so new shared strings land in the same order
# the original author produced them in. The row uses the grey "s=2"
# style throughout, so copy that formatting onto the new cell too
# (xlPasteFormats = -4122) before writing its text.
$ws.Range("D2").Value = "system_1 (-O3)"
$ws.Range("C2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "system_1 (-O0)"

# Re-type the header row for the columns that moved / appeared so the
# table keeps its simple numeric-looking labels "3".."11".
Set-TextValue $ws.Range("E1") "3"
Set-TextValue $ws.Range("F1") "4"
Set-TextValue $ws.Range("G1") "5"
Set-TextValue $ws.Range("H1") "6"
Set-TextValue $ws.Range("I1") "7"
Set-TextValue $ws.Range("J1") "8"
Set-TextValue $ws.Range("K1") "9"
Set-TextValue $ws.Range("L1") "10"
Set-TextValue $ws.Range("M1") "11"

# Row 3 (SW_version): the new configuration also uses SW version "1".
$ws.Range("E3").Value = 1

# Row 4 (no operation): updated baseline + two new measurements.
$ws.Range("C4").Value = 16
$ws.Range("D4").Value = 22
$ws.Range("E4").Value = 22

# Row 5 (simple add elapsed time(cycles)).
$ws.Range("D5").Value = 1.3
$ws.Range("E5").Value = 0

# Row 6 (simple mult elapsed time(cycles)).
$ws.Range("D6").Value = 1.5
$ws.Range("E6").Value = 0

# Row 7 (6x6 matrix mult elapsed time(cycles)).
$ws.Range("D7").Value = 1876
$ws.Range("E7").Value = 193

# Autofit the two new data columns, mirroring the width tweak Excel
# applies once real content lands in D:E.
$ws.Range("D1:E1").EntireColumn.AutoFit()

# Put the selection where the author left it.
$ws.Range("A3").Select()
